# Apply weekly fruit/vegetable price update by permuting row contents
# (columns D, M, N, O, P, R, S) across rows 2-15.
# Mapping: new row <- old row (source of the new values)
#   2 <- 12, 3 <- 13, 4 <- 10, 5 <- 5, 6 <- 6, 7 <- 4, 8 <- 2,
#   9 <- 9, 10 <- 14, 11 <- 11, 12 <- 15, 13 <- 8, 14 <- 7, 15 <- 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for the columns that change.
$cols = @("D", "M", "N", "O", "P", "R", "S")
$rows = 2..15

$original = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowData
}

# Mapping of destination row -> source row (values taken from source's
# original content and written into destination row).
$mapping = @{
    2  = 12
    3  = 13
    4  = 10
    5  = 5
    6  = 6
    7  = 4
    8  = 2
    9  = 9
    10 = 14
    11 = 11
    12 = 15
    13 = 8
    14 = 7
    15 = 3
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcData[$col]
    }
}
